# Apply "Atualização de bases das ligas" update to the Poland Ekstraklasa sheet.
#
# 1) Rename two header columns (shared strings used by I1/J1):
#       ht_goals_h -> HTHG
#       ht_goals_a -> HTAG
# 2) Rows 298-305 (match records, column A id stays put) get their data
#    (columns B..AD) cyclically re-shuffled among themselves, matching a
#    refreshed pull of the underlying match-odds feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header renames -----------------------------------------------------
$ws.Range("I1").Value2 = "HTHG"
$ws.Range("J1").Value2 = "HTAG"

# --- 2) Row data permutation for rows 298-305 ------------------------------
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R", `
          "S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

$rows = 298..305

# Read every source row's data (cols B..AD) into memory BEFORE writing
# anything, since the remap below is a permutation (not a simple shift) and
# source/destination ranges overlap.
$rowsData = @{}
foreach ($r in $rows) {
    $vals = @{}
    foreach ($col in $cols) {
        $vals[$col] = $ws.Range($col + $r).Value2
    }
    $rowsData[$r] = $vals
}

# Mapping: data currently sitting in row <key> is the data that should end
# up in row <value> (column A / the id column is untouched).
$rowMap = @{
    298 = 300
    299 = 301
    300 = 305
    301 = 302
    302 = 299
    303 = 304
    304 = 298
    305 = 303
}

foreach ($src in $rowMap.Keys) {
    $dst = $rowMap[$src]
    $vals = $rowsData[$src]
    foreach ($col in $cols) {
        $ws.Range($col + $dst).Value2 = $vals[$col]
    }
}
